$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new columns before column B, shifting existing data (B:BH) to (H:BN)
$ws.Range("B1:G1").EntireColumn.Insert()

# Match the header style used by the rest of row 1 (now starting at H1)
$ws.Range("H1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Populate the newly inserted header cells
$ws.Range("B1").Value = "Unnamed: 0.5"
$ws.Range("C1").Value = "Unnamed: 0.4"
$ws.Range("D1").Value = "Unnamed: 0.3"
$ws.Range("E1").Value = "Unnamed: 0.2"
$ws.Range("F1").Value = "Unnamed: 0.1"
$ws.Range("G1").Value = "Unnamed: 0"
